# Generate Report for Handoff
# Updates the localization-status report: a new handoff round has begun for
# both tracked files, so the "latest handoff" identifiers/timestamps move
# forward and the now-stale "latest handback" columns (F/G) are cleared out.

$wb = $excel.ActiveWorkbook

$srcMdA  = "e9f268f0-6fa4-405f-884e-dde653d1daaf.md"
$srcMdB  = "ffffcf185373-06c0-46a2-856d-8a53eb48effd.md"
$status  = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $srcMdA
$ws1.Range("B2").Value = $status
$ws1.Range("C2").Value = $status
$ws1.Range("D2").Value = "2016-51-17 18:51:54"

$ws1.Range("A3").Value = $srcMdB
$ws1.Range("B3").Value = $status
$ws1.Range("C3").Value = $status
$ws1.Range("D3").Value = "2016-51-17 18:51:54"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/12ca35d4-3a53-4132-b65d-5a124a67de2a.md", "", "", $srcMdA)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/e05a502e-0b09-4336-88b9-fa363d28dee0.md", "", "", $srcMdB)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$xlfA2 = "e9f268f0-6fa4-405f-884e-dde653d1daaf.ddc508a5c98821e1de20d56bd55a179aa83d7af8.zh-cn.xlf"

$ws2.Range("A2").Value = $srcMdA
$ws2.Range("C2").Value = $status
$ws2.Range("D2").Value = $xlfA2
$ws2.Range("E2").Value = "2016-03-17 18:51:50"
$ws2.Range("F2:G2").Clear()
$ws2.Range("H2").Value = "0001-01-01 00:00:00"

$ws2.Range("A3").Value = $srcMdB
$ws2.Range("C3").Value = $status
$ws2.Range("D3").Value = $xlfA2
$ws2.Range("E3").Value = "2016-03-17 18:51:50"
$ws2.Range("F3:G3").Clear()
$ws2.Range("H3").Value = "0001-01-01 00:00:00"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/12ca35d4-3a53-4132-b65d-5a124a67de2a.md", "", "", $srcMdA)
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/12ca35d4-3a53-4132-b65d-5a124a67de2a.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/913a889750a8e4e0fbcd5b38739a5194f543e3b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/12ca35d4-3a53-4132-b65d-5a124a67de2a.ed3b5d5dd956fa1588ae11ba18bc514cc6a12fd9.zh-cn.xlf", "", "", $xlfA2)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/65eee33ff430d897c0d83620cd5295b9075961dc/e2e/12ca35d4-3a53-4132-b65d-5a124a67de2a.md", "", "", $srcMdB)
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1af2cca5135da04040ea7d1beb539e514479397b/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/12ca35d4-3a53-4132-b65d-5a124a67de2a.ed3b5d5dd956fa1588ae11ba18bc514cc6a12fd9.zh-cn.xlf", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/e05a502e-0b09-4336-88b9-fa363d28dee0.md", "", "", $xlfA2)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$xlfA3 = "e9f268f0-6fa4-405f-884e-dde653d1daaf.ddc508a5c98821e1de20d56bd55a179aa83d7af8.de-de.xlf"

$ws3.Range("A2").Value = $srcMdA
$ws3.Range("C2").Value = $status
$ws3.Range("D2").Value = $xlfA3
$ws3.Range("E2").Value = "2016-03-17 18:51:54"
$ws3.Range("F2:G2").Clear()
$ws3.Range("H2").Value = "0001-01-01 00:00:00"

$ws3.Range("A3").Value = $srcMdB
$ws3.Range("C3").Value = $status
$ws3.Range("D3").Value = $xlfA3
$ws3.Range("E3").Value = "2016-03-17 18:51:54"
$ws3.Range("F3:G3").Clear()
$ws3.Range("H3").Value = "0001-01-01 00:00:00"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/12ca35d4-3a53-4132-b65d-5a124a67de2a.md", "", "", $srcMdA)
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/12ca35d4-3a53-4132-b65d-5a124a67de2a.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e85420b1dc98b61a4854be8f50dab573519a57b7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/12ca35d4-3a53-4132-b65d-5a124a67de2a.ed3b5d5dd956fa1588ae11ba18bc514cc6a12fd9.de-de.xlf", "", "", $xlfA3)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bd41380aa361878c957e8b2587a059194dd2f038/e2e/12ca35d4-3a53-4132-b65d-5a124a67de2a.md", "", "", $srcMdB)
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bb20b2653a129830c48f5084f8eed65e64ccdfab/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/12ca35d4-3a53-4132-b65d-5a124a67de2a.ed3b5d5dd956fa1588ae11ba18bc514cc6a12fd9.de-de.xlf", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/e05a502e-0b09-4336-88b9-fa363d28dee0.md", "", "", $xlfA3)
